# "OS Task Initial Release" — add two Change-Track log entries:
#   - a new row 2 ("OS Tick Initial Release", JMR, 00, Done)
#   - a new row 4 ("OS Task Initial Release", SPA, 02, On Process)
# while keeping the existing "Fix notification enable and disable." entry
# (old row 2) as row 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a fresh row above the existing data row, then stamp it with ---
# --- the same formatting (number formats / alignment) as that row so    ---
# --- the new cells don't inherit the header row's look.                 ---
$ws.Rows("2:2").Insert()
$ws.Range("A3:E3").Copy()
$ws.Range("A2:E2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill the new row 2 (order mirrors how the fields were actually typed)
$ws.Range("C2").Value = "JMR"
$ws.Range("D2").Value = "OS Tick Initial Release"
$ws.Range("B2").Value = "00"
$ws.Range("A2").Value = 41681
$ws.Range("E2").Value = "Done"

# --- Append a new row 4 with the same formatting as row 3 ---
$ws.Range("A3:E3").Copy()
$ws.Range("A4:E4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A4").Value = 41699
$ws.Range("B4").Value = "02"
$ws.Range("C4").Value = "SPA"
$ws.Range("D4").Value = "OS Task Initial Release"
$ws.Range("E4").Value = "On Process"

# Re-apply the "Change Number" column formatting (text, centered) — this
# also refreshes the header cell B1 to the combined centered/wrapped style.
$ws.Columns("B").NumberFormat = "@"
$ws.Columns("B").HorizontalAlignment = -4108

# Leave the selection where the editor ended up.
$ws.Range("D6").Select()
